# Updated cryptos list on Tue Nov 21 13:08:08 UTC 2023 with GitHub Actions
#
# Refreshes the Price (D) and Volume/1h (E) columns for each coin row.
# Some new Price values look like plain decimals (e.g. "259.74"); Excel's
# Range.Value setter would otherwise auto-convert those to numbers, which
# would change the cell's stored type away from the original inline text.
# A leading apostrophe forces text entry (exactly like typing it in the
# Excel UI); resetting Style back to "Normal" afterwards clears the
# quote-prefix/number-format style Excel applies so the cell's style index
# is left exactly as it was before the edit.
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$ws.Range("D2").Value = '37.237.91'
$ws.Range("E2").Value = '  +0.15%  '
$ws.Range("D3").Value = '2.011.47'
$ws.Range("E3").Value = '  -0.39%  '
$ws.Range("E4").Value = '  -0.03%  '
$ws.Range("D5").Value = '''259.74'
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '  +5.11%  '
$ws.Range("D6").Value = '''0.619'
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = '  -1.20%  '
$ws.Range("E7").Value = '  -0.05%  '
$ws.Range("D8").Value = '''56.59'
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = '  -6.48%  '
$ws.Range("E9").Value = '  -2.94%  '
$ws.Range("D10").Value = '''0.0771'
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = '  -5.18%  '
$ws.Range("E11").Value = '  -3.02%  '
$ws.Range("D12").Value = '''14.34'
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = '  -4.81%  '
$ws.Range("D13").Value = '2.306.32'
$ws.Range("E13").Value = '  -0.49%  '
$ws.Range("D14").Value = '''21.60'
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = '  -1.48%  '
$ws.Range("D15").Value = '''0.799'
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = '  -6.29%  '
$ws.Range("D16").Value = '''5.21'
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = '  -4.19%  '
$ws.Range("D17").Value = '2.034.03'
$ws.Range("E17").Value = '  +0.72%  '
$ws.Range("D18").Value = '37.124.34'
$ws.Range("E18").Value = '  -0.02%  '
$ws.Range("D19").Value = '''71.19'
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = '  +1.22%  '
$ws.Range("E20").Value = '  -3.76%  '
$ws.Range("D21").Value = '''233.82'
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = '  +1.45%  '
$ws.Range("D22").Value = '''5.10'
$ws.Range("D22").Style = "Normal"
$ws.Range("E23").Value = '  -0.13%  '
$ws.Range("D24").Value = '''2.59'
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = '  +0.36%  '
$ws.Range("E25").Value = '  -0.02%  '
$ws.Range("D26").Value = '''165.23'
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = '  +1.29%  '
$ws.Range("D27").Value = '''8.95'
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = '  -4.71%  '
$ws.Range("D28").Value = '''19.57'
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = '  -0.87%  '
$ws.Range("D29").Value = '''0.128'
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = '  -6.79%  '
$ws.Range("E30").Value = '  -4.04%  '
$ws.Range("E31").Value = '  -1.71%  '
$ws.Range("E32").Value = '  -3.78%  '
$ws.Range("D33").Value = '''0.0643'
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = '  -4.43%  '
$ws.Range("E34").Value = '  -0.43%  '
$ws.Range("E35").Value = '  -6.04%  '
$ws.Range("D36").Value = '''3.50'
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = '  -2.82%  '
$ws.Range("E37").Value = '  +0.60%  '
$ws.Range("E38").Value = '  +0.03%  '
$ws.Range("D39").Value = '''5.54'
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = '  +3.31%  '
$ws.Range("D40").Value = '''3.04'
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = '  +1.31%  '
$ws.Range("E41").Value = '  +0.17%  '
$ws.Range("D42").Value = '1.438.92'
$ws.Range("E42").Value = '  +4.43%  '
$ws.Range("D43").Value = '''0.0923'
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = '  -5.57%  '
$ws.Range("D44").Value = '''0.0210'
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = '  -2.07%  '
$ws.Range("D45").Value = '''89.39'
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = '  -2.32%  '
$ws.Range("D46").Value = '''15.63'
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = '  -7.64%  '
$ws.Range("E47").Value = '  -2.73%  '
$ws.Range("E48").Value = '  +2.02%  '
$ws.Range("D49").Value = '''6.96'
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = '  -6.40%  '
$ws.Range("D50").Value = '2.197.70'
$ws.Range("E51").Value = '  -7.81%  '
